$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge "Launch" + "Control" into a single cell D15
$ws.Range("D15").Value = "Launch Control"

# Remove B16 and C16 entirely (value + formula + formatting), D16 content cleared but keep the cell/style
$ws.Range("B16:C16").Clear()
$ws.Range("D16").ClearContents()

# Update the selection to match final state
[void]$ws.Range("E28").Select()
